$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87; existing rows 87:214 shift down to 88:215.
$ws.Rows("87").Insert()

# Populate the newly inserted row 87 with the new weekly record.
$ws.Cells.Item(87, 1).Value = 10
$ws.Cells.Item(87, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(87, 3).Value = "La Araucanía"
$ws.Cells.Item(87, 4).Value = 44495
$ws.Cells.Item(87, 5).Value = 9
$ws.Cells.Item(87, 6).Value = 100112044
$ws.Cells.Item(87, 7).Value = "Perejil"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 20
$ws.Cells.Item(87, 11).Value = 4000
$ws.Cells.Item(87, 12).Value = 4000
$ws.Cells.Item(87, 13).Value = 4000
$ws.Cells.Item(87, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(87, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(87, 16).Value = 1333
$ws.Cells.Item(87, 17).Value = 3
$ws.Cells.Item(87, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the other
# rows in column D (style carried over from the insert, but set explicitly
# to be safe).
$ws.Cells.Item(87, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
